$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resumo Obras")

# NORTE - Obras por FISCAL
$ws.Range("C15").Value = 6   # ISADORA ROSALINO
$ws.Range("C17").Value = 2   # LUCIANA POSTIÇO

# FISCAL / OBRAS summary
$ws.Range("C32").Value = 6   # ISADORA ROSALINO
$ws.Range("C35").Value = 17  # LUCIANA POSTIÇO
